# Swap the contents of columns G ("top_discipline") and H ("season"),
# including the header row, across the entire used range of the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $gCell = $ws.Cells.Item($r, 7)   # column G
    $hCell = $ws.Cells.Item($r, 8)   # column H

    $gValue = $gCell.Value()
    $hValue = $hCell.Value()

    $gCell.Value = $hValue
    $hCell.Value = $gValue
}
